# Daily attendance processing - normalize "Recorded By" (column G) entries.
# Rule observed in the source data: any literal "System" entries (exact case)
# are moved to the end of the comma-separated list (preserving the relative
# order of the remaining entries); if no "System" entry is present, the
# entries are sorted alphabetically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    $nonSystem = @()
    $systemParts = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $nonSystem += $p
        }
    }

    if ($systemParts.Count -gt 0) {
        $result = $nonSystem + $systemParts
    } else {
        $result = $parts | Sort-Object
    }

    $newVal = $result -join ", "

    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
    }
}
